# Update countries & provincias Spain
# Applies the diff: swap five pairs of country rows (name + re-ranked order)
# and refresh the numeric case counters, plus the "updated at" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 17:46"

# --- Swap country names whose ranking order changed ---
$ws.Range("A18").Value = "Irak"
$ws.Range("A19").Value = "Banglades"
$ws.Range("A37").Value = "Republica Dominicana"
$ws.Range("A38").Value = "Panama"
$ws.Range("A61").Value = "Moldavia"
$ws.Range("A62").Value = "Suiza"
$ws.Range("A100").Value = "Montenegro"
$ws.Range("A101").Value = "Eslovaquia"
$ws.Range("A183").Value = "Eritrea"
$ws.Range("A184").Value = "Mauricio"

# --- Refresh updated statistics values (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 7459200
$ws.Range("C4").Value = 11918
$ws.Range("D4").Value = 4712013
$ws.Range("E4").Value = 2535153
$ws.Range("G4").Value = 294
$ws.Range("H4").Value = 212034
$ws.Range("B17").Value = 460178
$ws.Range("C17").Value = 6914
$ws.Range("G17").Value = 59
$ws.Range("H17").Value = 42202
$ws.Range("B18").Value = 367474
$ws.Range("C18").Value = 4493
$ws.Range("D18").Value = 295882
$ws.Range("E18").Value = 62361
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = 9231
$ws.Range("B19").Value = 364987
$ws.Range("C19").Value = 1508
$ws.Range("D19").Value = 277078
$ws.Range("E19").Value = 82637
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 5272
$ws.Range("B22").Value = 317409
$ws.Range("C22").Value = 2548
$ws.Range("D22").Value = 228844
$ws.Range("E22").Value = 52647
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = 35918
$ws.Range("D25").Value = 257900
$ws.Range("E25").Value = 25708
$ws.Range("B29").Value = 160229
$ws.Range("C29").Value = 1471
$ws.Range("D29").Value = 136066
$ws.Range("E29").Value = 14847
$ws.Range("G29").Value = 19
$ws.Range("H29").Value = 9316
$ws.Range("B37").Value = 112728
$ws.Range("C37").Value = 519
$ws.Range("D37").Value = 88205
$ws.Range("E37").Value = 22415
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 2108
$ws.Range("B38").Value = 112595
$ws.Range("D38").Value = 89061
$ws.Range("E38").Value = 21162
$ws.Range("H38").Value = 2372
$ws.Range("B46").Value = 92409
$ws.Range("C46").Value = 663
$ws.Range("D46").Value = 80800
$ws.Range("E46").Value = 8348
$ws.Range("G46").Value = 15
$ws.Range("H46").Value = 3261
$ws.Range("B48").Value = 83563
$ws.Range("C48").Value = 553
$ws.Range("D48").Value = 76590
$ws.Range("E48").Value = 5402
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 1571
$ws.Range("B59").Value = 57784
$ws.Range("C59").Value = 19
$ws.Range("D59").Value = 57512
$ws.Range("E59").Value = 245
$ws.Range("B61").Value = 54064
$ws.Range("C61").Value = 1022
$ws.Range("D61").Value = 39499
$ws.Range("E61").Value = 13229
$ws.Range("G61").Value = 16
$ws.Range("H61").Value = 1336
$ws.Range("B62").Value = 53832
$ws.Range("C62").Value = 550
$ws.Range("D62").Value = 45300
$ws.Range("E62").Value = 6458
$ws.Range("H62").Value = 2074
$ws.Range("B92").Value = 14802
$ws.Range("C92").Value = 43
$ws.Range("D92").Value = 13961
$ws.Range("E92").Value = 508
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 333
$ws.Range("B95").Value = 13806
$ws.Range("C95").Value = 157
$ws.Range("D95").Value = 8077
$ws.Range("E95").Value = 5341
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 388
$ws.Range("B97").Value = 13101
$ws.Range("C97").Value = 1276
$ws.Range("D97").Value = 4752
$ws.Range("E97").Value = 8280
$ws.Range("G97").Value = 8
$ws.Range("H97").Value = 69
$ws.Range("B99").Value = 11373
$ws.Range("C99").Value = 108
$ws.Range("D99").Value = 9083
$ws.Range("E99").Value = 2167
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 123
$ws.Range("B100").Value = 10987
$ws.Range("C100").Value = 215
$ws.Range("D100").Value = 7397
$ws.Range("E100").Value = 3420
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 170
$ws.Range("B101").Value = 10938
$ws.Range("C101").Value = 797
$ws.Range("D101").Value = 4620
$ws.Range("E101").Value = 6270
$ws.Range("H101").Value = 48
$ws.Range("B116").Value = 6555
$ws.Range("C116").Value = 73
$ws.Range("D116").Value = 1991
$ws.Range("E116").Value = 4453
$ws.Range("G116").Value = 4
$ws.Range("H116").Value = 111
$ws.Range("B133").Value = 4542
$ws.Range("C133").Value = 11
$ws.Range("D133").Value = 2629
$ws.Range("E133").Value = 1837
$ws.Range("B140").Value = 3584
$ws.Range("C140").Value = 5
$ws.Range("D140").Value = 2216
$ws.Range("E140").Value = 1255
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 113
$ws.Range("B143").Value = 3382
$ws.Range("C143").Value = 2
$ws.Range("E143").Value = 136
$ws.Range("B156").Value = 2039
$ws.Range("C156").Value = 5
$ws.Range("D156").Value = 1297
$ws.Range("E156").Value = 155
$ws.Range("C183").Value = 6
$ws.Range("D183").Value = 353
$ws.Range("E183").Value = 28
$ws.Range("H183").Value = 0
$ws.Range("B184").Value = 381
$ws.Range("D184").Value = 344
$ws.Range("E184").Value = 27
$ws.Range("H184").Value = 10
